# Auto-generated edit script: updates crypto price/volume table cells
# to match the refreshed data from the GitHub Actions run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.935.86"
$ws.Range("E2").Value = "  -0.57%  "
$ws.Range("D3").Value = "2.329.56"
$ws.Range("E3").Value = "  +0.03%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "'302.86"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.24%  "
$ws.Range("D6").Value = "'95.06"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -2.86%  "
$ws.Range("E7").Value = "  -0.76%  "
$ws.Range("E8").Value = "  -0.07%  "
$ws.Range("D9").Value = "'0.496"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -1.62%  "
$ws.Range("D10").Value = "'34.10"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -4.34%  "
$ws.Range("D11").Value = "'19.02"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -1.35%  "
$ws.Range("E12").Value = "  -0.74%  "
$ws.Range("D13").Value = "'0.122"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +2.52%  "
$ws.Range("D14").Value = "'6.71"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -2.86%  "
$ws.Range("D15").Value = "2.692.88"
$ws.Range("E15").Value = "  +0.18%  "
$ws.Range("D16").Value = "2.319.91"
$ws.Range("E16").Value = "  -0.62%  "
$ws.Range("D17").Value = "'0.790"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +0.59%  "
$ws.Range("D18").Value = "42.882.13"
$ws.Range("E18").Value = "  -0.49%  "
$ws.Range("D19").Value = "'12.04"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -4.20%  "
$ws.Range("D20").Value = "'6.16"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +1.25%  "
$ws.Range("D21").Value = "0.0₃0891"
$ws.Range("E21").Value = "  -0.97%  "
$ws.Range("D22").Value = "'68.02"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +0.07%  "
$ws.Range("D23").Value = "'236.47"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -0.55%  "
$ws.Range("D24").Value = "'2.24"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +1.95%  "
$ws.Range("E25").Value = "  -0.13%  "
$ws.Range("E26").Value = "  -0.44%  "
$ws.Range("D27").Value = "'24.61"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -1.62%  "
$ws.Range("D28").Value = "'2.35"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -0.62%  "
$ws.Range("D29").Value = "'9.12"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +0.21%  "
$ws.Range("D30").Value = "'31.57"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -4.79%  "
$ws.Range("B31").Value = "Monero"
$ws.Range("C31").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D31").Value = "'140.87"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -15.04%  "
$ws.Range("B32").Value = "FirstDigitalUSD"
$ws.Range("C32").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D32").Value = "'1.00"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -0.03%  "
$ws.Range("E33").Value = "  -0.07%  "
$ws.Range("D34").Value = "'17.73"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -1.98%  "
$ws.Range("E35").Value = "  +1.18%  "
$ws.Range("E36").Value = "  -2.82%  "
$ws.Range("B37").Value = "ARBITRUM"
$ws.Range("C37").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D37").Value = "'1.80"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +2.29%  "
$ws.Range("B38").Value = "WEMIXToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D38").Value = "'2.30"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -1.64%  "
$ws.Range("E39").Value = "  -0.63%  "
$ws.Range("D40").Value = "'22.45"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +23.50%  "
$ws.Range("E41").Value = "  -1.43%  "
$ws.Range("E42").Value = "  -0.88%  "
$ws.Range("D43").Value = "1.935.65"
$ws.Range("E43").Value = "  -3.30%  "
$ws.Range("E44").Value = "  -0.78%  "
$ws.Range("D45").Value = "'10.06"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -5.29%  "
$ws.Range("E46").Value = "  -2.29%  "
$ws.Range("D47").Value = "'2.73"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -2.02%  "
$ws.Range("B48").Value = "HuobiToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D48").Value = "'2.87"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +1.68%  "
$ws.Range("B49").Value = "RocketPoolETH"
$ws.Range("C49").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D49").Value = "2.559.96"
$ws.Range("E49").Value = "  +0.21%  "
$ws.Range("D50").Value = "'53.19"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -0.99%  "
$ws.Range("D51").Value = "'72.54"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +0.60%  "
